$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the daily power record for row 75 (Date 2018-10-22 / serial 43399):
# Start Time = 0 (midnight), End Time = 0.34166666666666662 (~08:12:00 AM).
# The Duration / Second Duration / Absolute Value columns are driven by the
# existing shared formulas in D75:F75 and will recalculate automatically.
$ws.Range("B75").Value = 0
$ws.Range("C75").Value = 0.34166666666666662

# Update the active selection to match the edited workbook's view state.
$ws.Range("B72").Select()
